$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3.55
$ws.Range("F3").Value = 1.02
$ws.Range("H3").Value = 1.02
$ws.Range("F4").Value = 1.02
$ws.Range("H4").Value = 1.02
$ws.Range("F5").Value = 1.02
$ws.Range("H5").Value = 1.02
$ws.Range("R5").Value = 1.21
$ws.Range("F6").Value = 5.2
$ws.Range("G6").Value = 7.8
$ws.Range("H6").Value = 1.49
$ws.Range("I6").Value = 1.61
$ws.Range("J6").Value = 4.3
$ws.Range("K6").Value = 6.4
$ws.Range("N6").Value = 4.2
$ws.Range("P6").Value = 2.24
$ws.Range("Q6").Value = 1.61
$ws.Range("R6").Value = 1.49
$ws.Range("S6").Value = 2.56
$ws.Range("T6").Value = 1.78
$ws.Range("U6").Value = 1.98
$ws.Range("V6").Value = 2.62
$ws.Range("W6").Value = 1.16
$ws.Range("AF6").Value = 55
$ws.Range("AK6").Value = 85
$ws.Range("AL6").Value = 80
$ws.Range("AN6").Value = 90
$ws.Range("O7").Value = 1.16
$ws.Range("S7").Value = 2.2
$ws.Range("AB7").Value = 24
$ws.Range("AD7").Value = 12.5
$ws.Range("AE7").Value = 24
$ws.Range("AJ7").Value = 60
$ws.Range("AL7").Value = 40
$ws.Range("G8").Value = 2.38
$ws.Range("H8").Value = 3
$ws.Range("J8").Value = 3.75
$ws.Range("L8").Value = 1.01
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 2.4
$ws.Range("O8").Value = 1.18
$ws.Range("R8").Value = 1.49
$ws.Range("S8").Value = 2.18
$ws.Range("T8").Value = 1.04
$ws.Range("U8").Value = 1.04
$ws.Range("V8").Value = 1.34
$ws.Range("W8").Value = 1.72
$ws.Range("X8").Value = 34
$ws.Range("Y8").Value = 26
$ws.Range("Z8").Value = 38
$ws.Range("AA8").Value = 75
$ws.Range("AB8").Value = 21
$ws.Range("AC8").Value = 14.5
$ws.Range("AD8").Value = 21
$ws.Range("AE8").Value = 48
$ws.Range("AF8").Value = 25
$ws.Range("AG8").Value = 17
$ws.Range("AH8").Value = 22
$ws.Range("AI8").Value = 50
$ws.Range("AJ8").Value = 44
$ws.Range("AK8").Value = 30
$ws.Range("AL8").Value = 40
$ws.Range("AM8").Value = 85
$ws.Range("AN8").Value = 1000
$ws.Range("AO8").Value = 1000
$ws.Range("F9").Value = 5.7
$ws.Range("G9").Value = 6.8
$ws.Range("H9").Value = 1.6
$ws.Range("I9").Value = 1.73
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 4.6
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 4.1
$ws.Range("O9").Value = 1.28
$ws.Range("Q9").Value = 1.77
$ws.Range("R9").Value = 1.41
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 1.83
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 2.36
$ws.Range("W9").Value = 1.17
$ws.Range("X9").Value = 22
$ws.Range("Y9").Value = 990
$ws.Range("Z9").Value = 980
$ws.Range("AA9").Value = 980
$ws.Range("AB9").Value = 990
$ws.Range("AC9").Value = 990
$ws.Range("AD9").Value = 990
$ws.Range("AE9").Value = 980
$ws.Range("AF9").Value = 980
$ws.Range("AG9").Value = 990
$ws.Range("AH9").Value = 990
$ws.Range("AI9").Value = 980
$ws.Range("AJ9").Value = 180
$ws.Range("AK9").Value = 95
$ws.Range("AL9").Value = 90
$ws.Range("AM9").Value = 130
$ws.Range("AN9").Value = 110
$ws.Range("AO9").Value = 980
$ws.Range("G10").Value = 40
$ws.Range("H10").Value = 1.22
$ws.Range("I10").Value = 1.34
$ws.Range("J10").Value = 6.4
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 3.15
$ws.Range("O10").Value = 1.09
$ws.Range("P10").Value = 3.15
$ws.Range("R10").Value = 1.88
$ws.Range("S10").Value = 1.69
$ws.Range("T10").Value = 1.04
$ws.Range("U10").Value = 1.04
$ws.Range("V10").Value = 3.4
$ws.Range("W10").Value = 1.02
$ws.Range("X10").Value = 990
$ws.Range("Y10").Value = 990
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 990
$ws.Range("AC10").Value = 990
$ws.Range("AD10").Value = 990
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 990
$ws.Range("AH10").Value = 990
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000
$ws.Range("F11").Value = 4.3
$ws.Range("G11").Value = 7.8
$ws.Range("H11").Value = 1.58
$ws.Range("I11").Value = 1.9
$ws.Range("J11").Value = 3.65
$ws.Range("K11").Value = 7.6
$ws.Range("L11").Value = 1.01
$ws.Range("M11").Value = 1.01
$ws.Range("N11").Value = 2.22
$ws.Range("O11").Value = 1.19
$ws.Range("P11").Value = 2.2
$ws.Range("Q11").Value = 1.53
$ws.Range("R11").Value = 1.45
$ws.Range("S11").Value = 2.26
$ws.Range("T11").Value = 1.04
$ws.Range("U11").Value = 1.04
$ws.Range("V11").Value = 2.1
$ws.Range("W11").Value = 1.14
$ws.Range("X11").Value = 990
$ws.Range("Y11").Value = 990
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 990
$ws.Range("AC11").Value = 990
$ws.Range("AD11").Value = 990
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 1000
$ws.Range("AG11").Value = 990
$ws.Range("AH11").Value = 990
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("AO11").Value = 1000
$ws.Range("L12").Value = 1.01
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 2.36
$ws.Range("O12").Value = 1.15
$ws.Range("P12").Value = 2.34
$ws.Range("Q12").Value = 1.4
$ws.Range("R12").Value = 1.22
$ws.Range("S12").Value = 2.02
$ws.Range("T12").Value = 1.01
$ws.Range("U12").Value = 1.01
$ws.Range("V12").Value = 1.01
$ws.Range("W12").Value = 1.01
$ws.Range("X12").Value = 990
$ws.Range("Y12").Value = 990
$ws.Range("Z12").Value = 1000
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 990
$ws.Range("AC12").Value = 990
$ws.Range("AD12").Value = 990
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 1000
$ws.Range("AG12").Value = 990
$ws.Range("AH12").Value = 990
$ws.Range("AI12").Value = 1000
$ws.Range("AJ12").Value = 1000
$ws.Range("AK12").Value = 1000
$ws.Range("AL12").Value = 1000
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 1000
